# Append new Lancers job listings (scraped 2025-11-02 12:31:18 JST) to the
# "ランサーズ" sheet. The existing two rows (old r3 "Transformer..." and old
# r4 "ECフォース...") are pushed down to rows 7-8, and four brand-new rows
# are inserted at rows 3-6. All "取得日時" timestamps are refreshed to the
# new run time, and column B is widened to fit the longer titles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2025-11-02 12:31:18"

# --- widen the title column (B) -------------------------------------------
# ColumnWidth is in "characters"; Excel's internal px-rounding means we must
# feed it a value a little under 55 so the stored <col width="..."/> lands
# on an exact 55 (matches the target diff) instead of 55.8333....
$ws.Columns.Item(2).ColumnWidth = 54.14

# --- row 2: same listing, only the collection timestamp changes -----------
$ws.Range("A2").Value = $timestamp

# --- rows 3-6: brand-new listings ------------------------------------------
$ws.Range("A3").Value = $timestamp
$ws.Range("B3").Value = "【急募】AI×ノーコードで動画自動制作ワークフロー構築"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5425676"
$ws.Range("G3").Value = 303
$ws.Range("H3").Value = "🔥AI,Ai"

$ws.Range("A4").Value = $timestamp
$ws.Range("B4").Value = "製造業向けAI戦略アドバイザー募集(事業価値試算・プロジェクト推進支援)"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5419380"
$ws.Range("G4").Value = 298
$ws.Range("H4").Value = "🔥AI,Ai"

$ws.Range("A5").Value = $timestamp
$ws.Range("B5").Value = "初回 【急募】ECサイトの要件定義や基本設計ができる方を募集(1人月、フルリモート可、2025年12月〜)"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5425629"
$ws.Range("G5").Value = 45
$ws.Range("H5").Value = "◇サイト"

$ws.Range("A6").Value = $timestamp
$ws.Range("B6").Value = "対象ECサイトがどのECカートシステムを利用しているかの調査"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5425697"
$ws.Range("G6").Value = 45
$ws.Range("H6").Value = "◇サイト"

# --- rows 7-8: the two previously-existing listings, now pushed down ------
$ws.Range("A7").Value = $timestamp
$ws.Range("B7").Value = "【急募】Transformerベースのテキストエンコーダー経験者募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5425363"
$ws.Range("G7").Value = 25

$ws.Range("A8").Value = $timestamp
$ws.Range("B8").Value = "【急募】ECフォース EFO利用 META広告計測設定の経験者募集"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5425263"
$ws.Range("G8").Value = 13

# --- rebuild the hyperlinks on column F ------------------------------------
# This engine's Hyperlinks.Item(n).Delete() is a no-op, but Range.Hyperlinks
# .Delete() clears every hyperlink on the sheet in one go (regardless of the
# range it's called on); so wipe them all and re-add the seven links fresh,
# in row order, so the relationship ids come out rId1..rId7 matching F2..F8.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5405834")
$ws.Range("F2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5425676")
$ws.Range("F3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5419380")
$ws.Range("F4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5425629")
$ws.Range("F5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5425697")
$ws.Range("F6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5425363")
$ws.Range("F7").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5425263")
$ws.Range("F8").Style = "Hyperlink"
